$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default/no explicit format) taken from an untouched cell
# in column D so that re-applying it after a text-forcing write does not
# leave a numFmt/quotePrefix style on the edited cells.
$refStyle = $ws.Cells.Item(17, 4).Style

$ws.Cells.Item(2, 4).Value = "27.375.73"
$ws.Cells.Item(2, 5).Value = "  +4.96%  "
$ws.Cells.Item(3, 4).Value = "1.812.57"
$ws.Cells.Item(3, 5).Value = "  +5.52%  "
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 4).Style = $refStyle
$ws.Cells.Item(4, 5).Value = "  +0.26%  "
$ws.Cells.Item(5, 4).Value = "'342.34"
$ws.Cells.Item(5, 4).Style = $refStyle
$ws.Cells.Item(5, 5).Value = "  +2.90%  "
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 4).Style = $refStyle
$ws.Cells.Item(6, 5).Value = "  +0.14%  "
$ws.Cells.Item(7, 4).Value = "'0.3808"
$ws.Cells.Item(7, 4).Style = $refStyle
$ws.Cells.Item(7, 5).Value = "  +3.33%  "
$ws.Cells.Item(8, 4).Value = "'0.3496"
$ws.Cells.Item(8, 4).Style = $refStyle
$ws.Cells.Item(8, 5).Value = "  +4.43%  "
$ws.Cells.Item(9, 4).Value = "'48.79"
$ws.Cells.Item(9, 4).Style = $refStyle
$ws.Cells.Item(9, 5).Value = "  -1.02%  "
$ws.Cells.Item(10, 4).Value = "'1.232"
$ws.Cells.Item(10, 4).Style = $refStyle
$ws.Cells.Item(10, 5).Value = "  +3.94%  "
$ws.Cells.Item(11, 4).Value = "'0.07728"
$ws.Cells.Item(11, 4).Style = $refStyle
$ws.Cells.Item(11, 5).Value = "  +3.59%  "
$ws.Cells.Item(12, 4).Value = "'1.002"
$ws.Cells.Item(12, 4).Style = $refStyle
$ws.Cells.Item(12, 5).Value = "  +0.09%  "
$ws.Cells.Item(13, 2).Value = "Polkadot"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Cells.Item(13, 4).Value = "'6.630"
$ws.Cells.Item(13, 4).Style = $refStyle
$ws.Cells.Item(13, 5).Value = "  +5.38%  "
$ws.Cells.Item(14, 2).Value = "Solana"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(14, 4).Value = "'21.97"
$ws.Cells.Item(14, 4).Style = $refStyle
$ws.Cells.Item(14, 5).Value = "  +9.64%  "
$ws.Cells.Item(15, 4).Value = "'7.254"
$ws.Cells.Item(15, 4).Style = $refStyle
$ws.Cells.Item(15, 5).Value = "  +4.84%  "
$ws.Cells.Item(16, 4).Value = "1.809.33"
$ws.Cells.Item(16, 5).Value = "  +5.45%  "
$ws.Cells.Item(17, 5).Value = "  +3.87%  "
$ws.Cells.Item(18, 4).Value = "'0.06718"
$ws.Cells.Item(18, 4).Style = $refStyle
$ws.Cells.Item(18, 5).Value = "  +1.29%  "
$ws.Cells.Item(19, 4).Value = "'85.89"
$ws.Cells.Item(19, 4).Style = $refStyle
$ws.Cells.Item(19, 5).Value = "  +4.93%  "
$ws.Cells.Item(20, 4).Value = "'1.002"
$ws.Cells.Item(20, 4).Style = $refStyle
$ws.Cells.Item(20, 5).Value = "  +0.17%  "
$ws.Cells.Item(21, 4).Value = "'17.59"
$ws.Cells.Item(21, 4).Style = $refStyle
$ws.Cells.Item(21, 5).Value = "  +7.35%  "
$ws.Cells.Item(22, 4).Value = "'6.552"
$ws.Cells.Item(22, 4).Style = $refStyle
$ws.Cells.Item(22, 5).Value = "  +7.67%  "
$ws.Cells.Item(23, 4).Value = "'13.35"
$ws.Cells.Item(23, 4).Style = $refStyle
$ws.Cells.Item(23, 5).Value = "  +2.59%  "
$ws.Cells.Item(24, 4).Value = "27.417.12"
$ws.Cells.Item(24, 5).Value = "  +5.38%  "
$ws.Cells.Item(25, 4).Value = "'2.469"
$ws.Cells.Item(25, 4).Style = $refStyle
$ws.Cells.Item(25, 5).Value = "  -0.11%  "
$ws.Cells.Item(26, 4).Value = "'2.677"
$ws.Cells.Item(26, 4).Style = $refStyle
$ws.Cells.Item(26, 5).Value = "  +8.94%  "
$ws.Cells.Item(27, 4).Value = "'22.00"
$ws.Cells.Item(27, 4).Style = $refStyle
$ws.Cells.Item(27, 5).Value = "  +14.64%  "
$ws.Cells.Item(28, 4).Value = "'1.474"
$ws.Cells.Item(28, 4).Style = $refStyle
$ws.Cells.Item(28, 5).Value = "  +11.74%  "
$ws.Cells.Item(29, 4).Value = "'153.81"
$ws.Cells.Item(29, 4).Style = $refStyle
$ws.Cells.Item(29, 5).Value = "  +2.43%  "
$ws.Cells.Item(30, 4).Value = "2.009.40"
$ws.Cells.Item(30, 5).Value = "  +5.43%  "
$ws.Cells.Item(31, 4).Value = "'135.92"
$ws.Cells.Item(31, 4).Style = $refStyle
$ws.Cells.Item(31, 5).Value = "  +5.13%  "
$ws.Cells.Item(32, 4).Value = "'6.316"
$ws.Cells.Item(32, 4).Style = $refStyle
$ws.Cells.Item(32, 5).Value = "  +6.70%  "
$ws.Cells.Item(33, 4).Value = "'4.034"
$ws.Cells.Item(33, 4).Style = $refStyle
$ws.Cells.Item(33, 5).Value = "  -1.73%  "
$ws.Cells.Item(34, 5).Value = "  +7.78%  "
$ws.Cells.Item(35, 4).Value = "'0.08732"
$ws.Cells.Item(35, 4).Style = $refStyle
$ws.Cells.Item(35, 5).Value = "  +2.40%  "
$ws.Cells.Item(36, 4).Value = "'1.707"
$ws.Cells.Item(36, 4).Style = $refStyle
$ws.Cells.Item(36, 5).Value = "  -0.88%  "
$ws.Cells.Item(37, 4).Value = "'5.619"
$ws.Cells.Item(37, 4).Style = $refStyle
$ws.Cells.Item(37, 5).Value = "  +5.23%  "
$ws.Cells.Item(38, 4).Value = "'0.6987"
$ws.Cells.Item(38, 4).Style = $refStyle
$ws.Cells.Item(38, 5).Value = "  +13.37%  "
$ws.Cells.Item(39, 4).Value = "'0.2276"
$ws.Cells.Item(39, 4).Style = $refStyle
$ws.Cells.Item(39, 5).Value = "  +6.79%  "
$ws.Cells.Item(40, 4).Value = "'0.02418"
$ws.Cells.Item(40, 4).Style = $refStyle
$ws.Cells.Item(40, 5).Value = "  +5.35%  "
$ws.Cells.Item(41, 4).Value = "'0.06487"
$ws.Cells.Item(41, 4).Style = $refStyle
$ws.Cells.Item(41, 5).Value = "  +4.40%  "
$ws.Cells.Item(42, 4).Value = "'8.975"
$ws.Cells.Item(42, 4).Style = $refStyle
$ws.Cells.Item(42, 5).Value = "  +5.20%  "
$ws.Cells.Item(43, 5).Value = "  +6.50%  "
$ws.Cells.Item(44, 4).Value = "'14.69"
$ws.Cells.Item(44, 4).Style = $refStyle
$ws.Cells.Item(44, 5).Value = "  +1.77%  "
$ws.Cells.Item(45, 4).Value = "'0.6522"
$ws.Cells.Item(45, 4).Style = $refStyle
$ws.Cells.Item(45, 5).Value = "  +10.85%  "
$ws.Cells.Item(46, 5).Value = "  +0.14%  "
$ws.Cells.Item(47, 4).Value = "'4.025"
$ws.Cells.Item(47, 4).Style = $refStyle
$ws.Cells.Item(47, 5).Value = "  +4.99%  "
$ws.Cells.Item(48, 4).Value = "'2.184"
$ws.Cells.Item(48, 4).Style = $refStyle
$ws.Cells.Item(48, 5).Value = "  +8.24%  "
$ws.Cells.Item(49, 4).Value = "'132.72"
$ws.Cells.Item(49, 4).Style = $refStyle
$ws.Cells.Item(49, 5).Value = "  +3.90%  "
$ws.Cells.Item(50, 4).Value = "'0.07335"
$ws.Cells.Item(50, 4).Style = $refStyle
$ws.Cells.Item(50, 5).Value = "  +0.76%  "
$ws.Cells.Item(51, 4).Value = "'80.44"
$ws.Cells.Item(51, 4).Style = $refStyle
$ws.Cells.Item(51, 5).Value = "  +4.53%  "
